# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column C (municipio-nombre) and column E (aragon) metadata rows are
# updated to match the curated dimension pattern already used by columns
# D/H (refArea / dim / URI-*), and the now-unused "mapping-aragon.xlsx"
# value is cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

$ws.Range("C3").Value = "dim"

$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Comunidad"

$ws.Range("E5").Clear()
